$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the bookmark that currently sits at the end of the "OCR"
#    paragraph (it will be re-created in the new "Sponsor" paragraph).
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2) Merge the two runs "I" + "ntegrate an OCR service to consume
#    images" into a single run "Integrate an OCR service to consume
#    images" by doing a no-op Find & Replace over that exact text.
# ------------------------------------------------------------------
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("Integrate an OCR service to consume images", $true, $false, $false, $false, $false, $true, 1, $false, "Integrate an OCR service to consume images", 2)

# ------------------------------------------------------------------
# 3) First paragraph ("Project Name: ... (app name will be
#    different)") currently ends with two empty tab runs - drop them.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1End = $p1.Range.End
$tabsRange = $d.Range($p1End - 3, $p1End - 1)
$tabsRange.Delete()

# ------------------------------------------------------------------
# 4) Insert two new paragraphs right before the existing blank
#    paragraph that follows paragraph 1: one blank paragraph, and one
#    that will hold "Sponsor: Alea Bunker" (with the _GoBack bookmark
#    re-inserted between "Alea" and " Bunker").
# ------------------------------------------------------------------
$pBlankExisting = $d.Paragraphs(2)
$pBlankExisting.Range.InsertParagraphBefore()
$pBlankExisting = $d.Paragraphs(3)
$pBlankExisting.Range.InsertParagraphBefore()

$pSponsor = $d.Paragraphs(3)
$pSponsor.Range.InsertAfter("Sponsor: Alea Bunker")

# Re-insert the _GoBack bookmark right after "Sponsor: Alea" (position
# 13 characters into the new paragraph), splitting " Bunker" into its
# own run, matching the target markup.
$pSponsor = $d.Paragraphs(3)
$bmPos = $pSponsor.Range.Start + "Sponsor: Alea".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
